$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 2
$ws.Range("A3").Formula = "=AVERAGE(A1:A2)"
$ws.Range("A4").Select()
